$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E25").Value = 2.000000000000002
$ws.Range("F25").Value = 2.051282051282053
$ws.Range("K25").Value = 50.78544871794873
$ws.Range("F26").Value = 16
$ws.Range("K26").Value = 63.565
$ws.Range("H95").Value = 48.77749999999762
$ws.Range("J95").Value = 8.572500000002385
$ws.Range("K95").Value = 48.12291666666428
$ws.Range("H96").Value = 48.77749999999762
$ws.Range("H97").Value = 48.77749999999762
$ws.Range("H98").Value = 48.77749999999762
$ws.Range("H99").Value = 48.77749999999762
$ws.Range("H100").Value = 76.62749999999762
$ws.Range("H101").Value = 104.7274999999976
$ws.Range("I101").Value = 28.1
$ws.Range("K101").Value = 76.32416666666667
$ws.Range("H102").Value = 104.7274999999976
$ws.Range("H103").Value = 104.7274999999976
$ws.Range("H104").Value = 104.4774999999976
$ws.Range("H105").Value = 104.1524999999976
$ws.Range("H106").Value = 102.9024999999976
$ws.Range("H107").Value = 101.2274999999976
$ws.Range("H108").Value = 99.20249999999761
$ws.Range("E109").Value = 34.44102564102565
$ws.Range("G109").Value = 12.42
$ws.Range("H109").Value = 96.82749999999761
$ws.Range("K109").Value = 18.42583333333334
$ws.Range("E110").Value = 34.44102564102565
$ws.Range("H110").Value = 93.87749999999761
$ws.Range("E111").Value = 50.04102564102565
$ws.Range("F111").Value = 16
$ws.Range("H111").Value = 90.8524999999976
$ws.Range("K111").Value = 43.75708333333333
$ws.Range("E112").Value = 65.64102564102565
$ws.Range("F112").Value = 16
$ws.Range("H112").Value = 86.92749999999761
$ws.Range("K112").Value = 51.6675
$ws.Range("E113").Value = 49.23076923076924
$ws.Range("H113").Value = 81.75249999999761
$ws.Range("E114").Value = 32.82051282051282
$ws.Range("H114").Value = 72.97749999999762
$ws.Range("F115").Value = 0
$ws.Range("G115").Value = 16
$ws.Range("H115").Value = 64.27749999999762
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 8.699999999999999
$ws.Range("K115").Value = 40.54958333333334
$ws.Range("H116").Value = 56.37749999999762
$ws.Range("H117").Value = 48.35249999999762
$ws.Range("H118").Value = 39.17749999999761
$ws.Range("H119").Value = 30.40249999999762
$ws.Range("H120").Value = 30.40249999999762
$ws.Range("H121").Value = 30.40249999999762
$ws.Range("H122").Value = 30.40249999999762
$ws.Range("H123").Value = 30.40249999999762
$ws.Range("H124").Value = 47.95
$ws.Range("I124").Value = 17.54750000000239
$ws.Range("K124").Value = 63.01416666666906
$ws.Range("H125").Value = 47.95
$ws.Range("H126").Value = 47.95
$ws.Range("H127").Value = 47.95
$ws.Range("H128").Value = 47.95
$ws.Range("H129").Value = 47.95
$ws.Range("H130").Value = 47.95
$ws.Range("H131").Value = 47.95
$ws.Range("H132").Value = 46.6
$ws.Range("H133").Value = 44.575
$ws.Range("H134").Value = 41.825
$ws.Range("H135").Value = 37.85
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 41.61666666666667
$ws.Range("H366").Value = 65.3
$ws.Range("I366").Value = 7.724999999999994
$ws.Range("K366").Value = 64.33291666666666
$ws.Range("H367").Value = 65.3
$ws.Range("J368").Value = 0
$ws.Range("K368").Value = 46.6625
$ws.Range("H435").Value = 1.199999999999985
$ws.Range("I435").Value = 1.199999999999985
$ws.Range("K435").Value = 45.70583333333331
$ws.Range("H436").Value = 29.04999999999999
$ws.Range("H437").Value = 57.14999999999998
$ws.Range("H438").Value = 85.17499999999998
$ws.Range("H439").Value = 85.17499999999998
$ws.Range("H440").Value = 85.17499999999998
$ws.Range("H441").Value = 85.17499999999998
$ws.Range("J442").Value = 0
$ws.Range("K442").Value = 45.81541666666667
$ws.Range("H651").Value = 24.2
$ws.Range("I651").Value = 0
$ws.Range("K651").Value = 45.38083333333333
$ws.Range("H652").Value = 24.2
$ws.Range("H653").Value = 24.2
$ws.Range("H654").Value = 24.2
$ws.Range("H655").Value = 24.2
$ws.Range("H656").Value = 23.45
$ws.Range("H657").Value = 22.95
$ws.Range("H658").Value = 22.125
$ws.Range("H659").Value = 20.875
$ws.Range("H660").Value = 19.525
$ws.Range("H661").Value = 17.5
$ws.Range("H662").Value = 14.75
$ws.Range("H663").Value = 14.75
$ws.Range("H664").Value = 14.75
$ws.Range("H665").Value = 14.75
$ws.Range("H666").Value = 8.300000000000001
$ws.Range("H667").Value = 0.8750000000000009
$ws.Range("J668").Value = 0.8750000000000009
$ws.Range("K668").Value = 55.41500000000001
$ws.Range("H674").Value = 53.625
$ws.Range("I674").Value = 27.35
$ws.Range("K674").Value = 70.09708333333333
$ws.Range("I675").Value = 23.625
$ws.Range("K675").Value = 68.53083333333333
